$wb = $excel.ActiveWorkbook

# --- Moorings sheet: update glider reference designator from GL003 to GL365 ---
$moorings = $wb.Worksheets.Item("Moorings")
$moorings.Activate()
$moorings.Range("A2").Value = "GP05MOAS-GL365"

# Move the active selection on the Moorings sheet to D10 (as reflected in the saved view)
$moorings.Range("D10").Select()

# --- Asset_Cal_Info sheet: update all instrument reference designators from GL003 to GL365 ---
$calInfo = $wb.Worksheets.Item("Asset_Cal_Info")

$calInfo.Range("A3").Value = "GP05MOAS-GL365-00-ENG000000"
$calInfo.Range("A4").Value = "GP05MOAS-GL365-01-FLORDM000"
$calInfo.Range("A5").Value = "GP05MOAS-GL365-01-FLORDM000"
$calInfo.Range("A6").Value = "GP05MOAS-GL365-01-FLORDM000"
$calInfo.Range("A7").Value = "GP05MOAS-GL365-01-FLORDM000"
$calInfo.Range("A8").Value = "GP05MOAS-GL365-02-DOSTAM000"
$calInfo.Range("A9").Value = "GP05MOAS-GL365-04-CTDGVM000"
